# Apply odds updates scraped for Jogos_da_Semana_FlashScore_2025-06-04.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 3 (Valour vs Vancouver FC) ----
$ws.Range("I3").Value  = 2.82
$ws.Range("K3").Value  = 9.25
$ws.Range("L3").Value  = 1.21
$ws.Range("M3").Value  = 4.15
$ws.Range("N3").Value  = 1.62
$ws.Range("O3").Value  = 2.2
$ws.Range("P3").Value  = 1.32
$ws.Range("Q3").Value  = 3.2
$ws.Range("R3").Value  = 1.53
$ws.Range("S3").Value  = 2.35
$ws.Range("T3").Value  = 10
$ws.Range("U3").Value  = 14
$ws.Range("V3").Value  = 9.5
$ws.Range("W3").Value  = 25
$ws.Range("X3").Value  = 17.5
$ws.Range("Y3").Value  = 23
$ws.Range("Z3").Value  = 9.25
$ws.Range("AA3").Value = 7.6
$ws.Range("AC3").Value = 45
$ws.Range("AD3").Value = 300
$ws.Range("AE3").Value = 11.25
$ws.Range("AJ3").Value = 26

# ---- Row 5 (America De Cali vs Junior) ----
$ws.Range("G5").Value = 2
$ws.Range("I5").Value = 4.2

# ---- Row 12 (Fenix vs Uruguay Montevideo) ----
$ws.Range("G12").Value  = 2.12
$ws.Range("I12").Value  = 3.45
$ws.Range("L12").Value  = 1.52
$ws.Range("M12").Value  = 2.22
$ws.Range("N12").Value  = 2.47
$ws.Range("P12").Value  = 1.55
$ws.Range("Q12").Value  = 2.15
$ws.Range("R12").Value  = 2.15
$ws.Range("S12").Value  = 1.55
$ws.Range("T12").Value  = 5.5
$ws.Range("U12").Value  = 8.75
$ws.Range("W12").Value  = 20
$ws.Range("X12").Value  = 22
$ws.Range("Z12").Value  = 6.2
$ws.Range("AA12").Value = 6.1
$ws.Range("AB12").Value = 20
$ws.Range("AE12").Value = 7.4
$ws.Range("AF12").Value = 16
$ws.Range("AG12").Value = 13.5
$ws.Range("AH12").Value = 50
$ws.Range("AJ12").Value = 65
